# daily auto push: 2026-02-27 19:00 UTC
#
# Inserts one new row of data (2026/02/28, Saturday, hour 1) into the
# "sei2" time-series sheet at row 889, pushing the existing 2026/12/29 ...
# 2027/01/05 rows (old rows 889-930) down by one (new rows 890-931).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at position 889; everything below (old rows
# 889..930, the 2026/12/29 .. 2027/01/05 block) shifts down to 890..931.
$ws.Rows.Item(889).Insert()

# Column A holds dates as plain text (e.g. "2026/12/29"), not real Excel
# date serials, everywhere else in the sheet. Force the cell to Text
# before assigning so the "YYYY/MM/DD"-looking string isn't silently
# reinterpreted as a date value, then drop the format override again so
# the cell ends up with the same (default/no-style) formatting as its
# neighbours.
$ws.Cells.Item(889, 1).NumberFormat = "@"
$ws.Cells.Item(889, 1).Value = "2026/02/28"
$ws.Cells.Item(889, 1).ClearFormats()

$ws.Cells.Item(889, 2).Value = "土"
$ws.Cells.Item(889, 3).Value = 1
$ws.Cells.Item(889, 4).Value = 201
